# Daily attendance processing - 2025-12-02 09:32:42
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session's attendance, separated by ", ". Entries that were recorded
# by "System" followed by a single real user account need to have that
# ordering reversed, e.g.:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#   "System, admin@admin.com"     ->  "admin@admin.com, System"
#
# Rows whose "System, ..." entry refers to the backup account
# (backup@backdoor.com) are left untouched, as are rows that already
# list the user before "System", rows with more than two recorders,
# and rows with only "System" and no other recorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $v = $cell.Text

    if ($v -like "System, *" -and $v -notlike "*backup@backdoor.com*") {
        $rest = $v.Substring(8)   # text after "System, "
        if ($rest -notlike "*,*") {
            $cell.Value = $rest + ", System"
        }
    }
}
